$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a date serial (formatted via numFmtId 14). Move it forward one month.
$ws.Range("A1").Value = 45436

# Update the price list (column D, rows 33-38) with the new amounts.
$ws.Range("D33").Value = 203.074
$ws.Range("D34").Value = 162.027
$ws.Range("D35").Value = 151.226
$ws.Range("D36").Value = 347.818
$ws.Range("D37").Value = 248.442
$ws.Range("D38").Value = 218.196
